$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new "Price" values look like plain decimal numbers
# (e.g. "580.79"). The source workbook stores every Price/Volume cell as
# text (inline string), so for those cells we briefly switch the cell to a
# text number format before assigning the value (this stops Excel from
# auto-converting the text into a floating point number), then clear the
# format again so the cell is left with no extra style, matching the rest
# of the sheet.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D5").Value = "580.79"
$ws.Range("D6").Value = "175.82"
$ws.Range("D8").Value = "0.588"
$ws.Range("D11").Value = "0.577"
$ws.Range("D12").Value = "45.32"
$ws.Range("D14").Value = "669.44"
$ws.Range("D20").Value = "17.38"
$ws.Range("D21").Value = "10.96"
$ws.Range("D23").Value = "5.41"
$ws.Range("D24").Value = "17.08"
$ws.Range("D25").Value = "98.91"
$ws.Range("D29").Value = "33.57"
$ws.Range("D30").Value = "8.42"
$ws.Range("D31").Value = "7.31"
$ws.Range("D32").Value = "571.16"
$ws.Range("D35").Value = "0.999"
$ws.Range("D38").Value = "3.32"
$ws.Range("D39").Value = "34.27"
$ws.Range("D43").Value = "3.31"
$ws.Range("D51").Value = "129.26"

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D51").ClearFormats()

# Remaining cells already stay text on their own (they contain characters
# Excel will not parse as a plain number, e.g. extra "." separators,
# a leading/trailing space, a "%" sign, or a subscript digit).
$ws.Range("D2").Value = "67.705.49"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.331.72"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("E6").Value = "  -3.25%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "3.327.07"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("E14").Value = "  +5.29%  "
$ws.Range("D15").Value = "3.877.17"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "67.735.88"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "3.336.84"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("E23").Value = "  +8.07%  "
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("E25").Value = "  +1.94%  "
$ws.Range("E26").Value = "  -3.57%  "
$ws.Range("E27").Value = "  -3.47%  "
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("E31").Value = "  +10.53%  "
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "3.686.96"
$ws.Range("E36").Value = "  -6.04%  "
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("E38").Value = "  -5.42%  "
$ws.Range("E39").Value = "  +5.58%  "
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("E42").Value = "  -4.19%  "
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("D45").Value = "0.0₃0665"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("E51").Value = "  -0.48%  "
